$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting from column J (2022) into new column K (2023) by copying
# the formats (number format, borders, fill, font, alignment) of the existing
# last column so the new column matches the established table style.
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new 2023 column with its values.
$ws.Cells.Item(3, 11).Value2 = 2023
$ws.Cells.Item(4, 11).Value2 = 1580.7
$ws.Cells.Item(5, 11).Value2 = 1193.3
$ws.Cells.Item(6, 11).Value2 = 1867.1
